{"js": "// The paragraph currently reads (across many runs, split on spaces):\n//   \"Subitamente, da Caixa de M\u00fasica Sint\u00e9tica, uma Voz come\u00e7ou a falar. A Voz da\n//    Raz\u00e3o, a Voz da Benevol\u00eancia. O cilindro girava com o Discurso Sint\u00e9tico N\u00famero Dois\n//    (For\u00e7a M\u00e9dia) Contra Motins, ... Pe\u00e7o-lhes, por favor, sejam bons e...\"\"\n//\n// The edit trims the paragraph down to just the opening sentence(s):\n//   \"Subitamente, da Caixa de M\u00fasica Sint\u00e9tica, uma Voz come\u00e7ou a falar. A Voz da\n//    Raz\u00e3o, a Voz da Benevol\u00eancia. \"\n// while keeping the _GoBack bookmark that sits right before the paragraph end.\n\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst para = paragraphs.items[0];\n\nconst keepText =\n  \"Subitamente, da Caixa de M\u00fasica Sint\u00e9tica, uma Voz come\u00e7ou a falar. \" +\n  \"A Voz da Raz\u00e3o, a Voz da Benevol\u00eancia. \";\n\n// 1) Find where the text to discard begins, and delete everything from\n//    there through the end of the paragraph. This leaves the bookmark\n//    (which sits right at the end of the paragraph, after all the runs)\n//    untouched.\nconst cutStart = body.search(\"O cilindro girava\", { matchCase: true });\ncutStart.load(\"items\");\nawait context.sync();\n\nconst paragraphRange = para.getRange();\nconst paragraphEnd = paragraphRange.getRange(\"End\");\nconst cutRange = cutStart.items[0].expandTo(paragraphEnd);\ncutRange.insertText(\"\", \"Replace\");\nawait context.sync();\n\n// 2) Re-write the remaining (now shorter) paragraph range with the target\n//    text so the surviving runs collapse into a single run, matching the\n//    merged <w:r> produced by the real edit.\nconst finalRange = para.getRange();\nfinalRange.insertText(keepText, \"Replace\");\nawait context.sync();\n", "ps1": "# The paragraph currently reads (split across many runs, one per word-wrap\n# boundary, with single-space runs in between):\n#   \"Subitamente, da Caixa de M\u00fasica Sint\u00e9tica, uma Voz come\u00e7ou a falar. A Voz da\n#    Raz\u00e3o, a Voz da Benevol\u00eancia. O cilindro girava com o Discurso Sint\u00e9tico N\u00famero Dois\n#    (For\u00e7a M\u00e9dia) Contra Motins, ... Pe\u00e7o-lhes, por favor, sejam bons e...\"\"\n#\n# The edit trims the paragraph down to just the opening sentence(s):\n#   \"Subitamente, da Caixa de M\u00fasica Sint\u00e9tica, uma Voz come\u00e7ou a falar. A Voz da\n#    Raz\u00e3o, a Voz da Benevol\u00eancia. \"\n# while keeping the _GoBack bookmark that sits right before the final run\n# (\"desejo que voc\u00eas sejam bons! ...\"). Editing a Range that straddles a\n# bookmark's (zero-width) anchor point drops the bookmark, so the deletions\n# below are done in two passes that each stay entirely on one side of it.\n\n$d = $word.ActiveDocument\n\n$finalText = \"Subitamente, da Caixa de M\u00fasica Sint\u00e9tica, uma Voz come\u00e7ou a falar. A Voz da Raz\u00e3o, a Voz da Benevol\u00eancia. \"\n\n# 1) Remove the very last run's text (\"desejo que voc\u00eas sejam bons! ...\") --\n#    this portion sits entirely after the bookmark, so deleting it does not\n#    cross the bookmark's anchor.\n$rng = $d.Content\n$rng.Find.Execute(\"desejo que voc\u00eas sejam bons! Pe\u00e7o-lhes, por favor, sejam bons e...\") | Out-Null\n$tail = $d.Range($rng.Start, $d.Content.End)\n$tail.Text = \"\"\n\n# 2) Remove everything from \"O cilindro girava\" through the (now-adjacent)\n#    end of content -- this portion sits entirely before the bookmark, so it\n#    likewise never crosses the anchor point.\n$rng2 = $d.Content\n$rng2.Find.Execute(\"O cilindro girava\") | Out-Null\n$mid = $d.Range($rng2.Start, $d.Content.End)\n$mid.Text = \"\"\n\n# 3) Rewrite the remaining (now short) text before the bookmark so the\n#    surviving runs collapse into a single run, matching the merged <w:r>\n#    produced by the real edit.\n$full = $d.Range(0, $d.Content.End)\n$full.Text = $finalText\n"}
